$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = 1
$ws.Range("D8").Value = ""
$ws.Rows.Item(8).RowHeight = 28

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = ""

$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "Fbeta avec beta = 2"

$ws.Range("C14").Value = 0.9
$ws.Range("D14").Value = "Fait dans une pipeline"

$ws.Range("C20").Value = 0.8
$ws.Range("D20").Value = "GridSearch pour l'optimisation"

$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Tableau à faire"

$ws.Range("C21").Value = 0.7
$ws.Range("D21").Value = "SHAP / LIME "

$ws.Range("C24").Value = 1
$ws.Range("D24").Value = ""

$ws.Range("C25").Value = 1

$ws.Range("C26").Value = 1

$ws.Range("D12").Select()
